# Auto-generated edit script applying cryptos.xlsx row updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.815.98'
$ws.Range('E2').Value = '  -0.11%  '
$ws.Range('D3').Value = '1.642.92'
$ws.Range('E3').Value = '  -0.30%  '
$ws.Range('E4').Value = '  -0.41%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '218.36'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.50%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.501'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.11%  '
$ws.Range('E7').Value = '  -0.46%  '
$ws.Range('E8').Value = '  -0.32%  '
$ws.Range('E9').Value = '  -1.02%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.22'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.03%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0849'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.64%  '
$ws.Range('D12').Value = '1.871.70'
$ws.Range('E12').Value = '  -0.22%  '
$ws.Range('D13').Value = '1.639.67'
$ws.Range('E13').Value = '  -0.65%  '
$ws.Range('E14').Value = '  -0.59%  '
$ws.Range('E15').Value = '  -0.68%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.18'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.91%  '
$ws.Range('D17').Value = '26.820.23'
$ws.Range('E17').Value = '  -0.09%  '
$ws.Range('E18').Value = '  -0.62%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '215.54'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.62%  '
$ws.Range('E20').Value = '  -0.42%  '
$ws.Range('E21').Value = '  +5.41%  '
$ws.Range('E22').Value = '  -0.13%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.35'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.14%  '
$ws.Range('E24').Value = '  -2.09%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '147.67'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.65%  '
$ws.Range('E26').Value = '  -0.68%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.118'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.10%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.13'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.29%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.73'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.02%  '
$ws.Range('E30').Value = '  -0.87%  '
$ws.Range('E31').Value = '  +0.94%  '
$ws.Range('E32').Value = '  +1.78%  '
$ws.Range('E33').Value = '  -0.74%  '
$ws.Range('E34').Value = '  +0.77%  '
$ws.Range('D35').Value = '1.268.78'
$ws.Range('E35').Value = '  -1.69%  '
$ws.Range('E36').Value = '  +0.06%  '
$ws.Range('E37').Value = '  +0.43%  '
$ws.Range('E38').Value = '  -1.42%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.817'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.42%  '
$ws.Range('E40').Value = '  -0.53%  '
$ws.Range('E41').Value = '  -0.94%  '
$ws.Range('E42').Value = '  -0.54%  '
$ws.Range('D43').Value = '1.781.30'
$ws.Range('E43').Value = '  -0.82%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.13'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.78%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '92.77'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.28%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '61.37'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.58%  '
$ws.Range('E47').Value = '  -0.22%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0515'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.51%  '
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0964'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.50%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.53'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.24%  '
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.01'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.42%  '
